$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number & report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/27/2023  Through  12/3/2023"

# --- Crime data table updates (rows 15-30) ---
$ws.Range("F15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1

$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100

$ws.Range("F15").Value = 2

$ws.Range("G15").Value = 2

$ws.Range("H15").Value = 0

$ws.Range("I15").Value = 20

$ws.Range("J15").Value = 20

$ws.Range("K15").Value = 0

$ws.Range("L15").Value = 122.222222222222

$ws.Range("M15").Value = 81.818181818181

$ws.Range("N15").Value = 5.263157894736

$ws.Range("D16").Value = 3

$ws.Range("E16").Value = -66.666666666666

$ws.Range("F16").Value = 7

$ws.Range("G16").Value = 10

$ws.Range("H16").Value = -30

$ws.Range("I16").Value = 86

$ws.Range("J16").Value = 89

$ws.Range("K16").Value = -3.370786516853

$ws.Range("L16").Value = 16.216216216216

$ws.Range("M16").Value = -49.411764705882

$ws.Range("N16").Value = -87.426900584795

$ws.Range("C17").Value = 2

$ws.Range("E17").Value = -33.333333333333

$ws.Range("F17").Value = 14

$ws.Range("H17").Value = 27.272727272727

$ws.Range("I17").Value = 204

$ws.Range("J17").Value = 199

$ws.Range("K17").Value = 2.51256281407

$ws.Range("L17").Value = 33.333333333333

$ws.Range("M17").Value = 31.612903225806

$ws.Range("N17").Value = -31.313131313131

$ws.Range("C18").Value = 5

$ws.Range("D18").Value = 2

$ws.Range("E18").Value = 150

$ws.Range("F18").Value = 8

$ws.Range("H18").Value = -46.666666666666

$ws.Range("I18").Value = 107

$ws.Range("J18").Value = 173

$ws.Range("K18").Value = -38.150289017341

$ws.Range("L18").Value = -29.605263157894

$ws.Range("M18").Value = -70.194986072423

$ws.Range("N18").Value = -93.206349206349

$ws.Range("C19").Value = 9

$ws.Range("D19").Value = 11

$ws.Range("E19").Value = -18.181818181818

$ws.Range("F19").Value = 29

$ws.Range("G19").Value = 48

$ws.Range("H19").Value = -39.583333333333

$ws.Range("I19").Value = 512

$ws.Range("J19").Value = 588

$ws.Range("K19").Value = -12.925170068027

$ws.Range("L19").Value = 20.18779342723

$ws.Range("M19").Value = 44.225352112676

$ws.Range("N19").Value = -17.419354838709

$ws.Range("C20").Value = 3

$ws.Range("D20").Value = 2

$ws.Range("E20").Value = 50

$ws.Range("F20").Value = 8

$ws.Range("G20").Value = 7

$ws.Range("H20").Value = 14.285714285714

$ws.Range("I20").Value = 145

$ws.Range("J20").Value = 110

$ws.Range("K20").Value = 31.818181818181

$ws.Range("L20").Value = 113.235294117647

$ws.Range("M20").Value = 15.079365079365

$ws.Range("N20").Value = -91.425192193968

$ws.Range("C21").Value = 20

$ws.Range("D21").Value = 22

$ws.Range("E21").Value = -9.090909090909

$ws.Range("F21").Value = 68

$ws.Range("G21").Value = 93

$ws.Range("H21").Value = -26.881720430107

$ws.Range("I21").Value = 1076

$ws.Range("J21").Value = 1180

$ws.Range("K21").Value = -8.813559322033

$ws.Range("L21").Value = 21.857304643261

$ws.Range("M21").Value = -8.503401360544

$ws.Range("N21").Value = -78.022875816993

$ws.Range("F22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1

$ws.Range("K22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100

$ws.Range("F22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1

$ws.Range("K22").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("H22").Value = 0

$ws.Range("J22").Value = 5

$ws.Range("K22").Value = 80

$ws.Range("C24").Value = 13

$ws.Range("D24").Value = 24

$ws.Range("E24").Value = -45.833333333333

$ws.Range("F24").Value = 70

$ws.Range("G24").Value = 83

$ws.Range("H24").Value = -15.662650602409

$ws.Range("I24").Value = 984

$ws.Range("J24").Value = 1035

$ws.Range("K24").Value = -4.927536231884

$ws.Range("L24").Value = 21.182266009852

$ws.Range("M24").Value = 16.725978647686

$ws.Range("C25").Value = 3

$ws.Range("D25").Value = 5

$ws.Range("E25").Value = -40

$ws.Range("F25").Value = 17

$ws.Range("G25").Value = 30

$ws.Range("H25").Value = -43.333333333333

$ws.Range("I25").Value = 323

$ws.Range("J25").Value = 327

$ws.Range("K25").Value = -1.223241590214

$ws.Range("L25").Value = 10.616438356164

$ws.Range("M25").Value = -15.22309711286

$ws.Range("F26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 2

$ws.Range("H26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100

$ws.Range("F26").Value = 5

$ws.Range("G26").Value = 4

$ws.Range("H26").Value = 25

$ws.Range("I26").Value = 25

$ws.Range("J26").Value = 29

$ws.Range("K26").Value = -13.793103448275

$ws.Range("L26").Value = 78.571428571428

$ws.Range("D27").Value = 1

$ws.Range("F27").Value = 3

$ws.Range("H27").Value = -62.5

$ws.Range("J27").Value = 76

$ws.Range("K27").Value = -19.736842105263

$ws.Range("L27").Value = -7.575757575757

$ws.Range("C28").Formula = "=""0"""
$ws.Range("C28").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("C29").Formula = "=""0"""
$ws.Range("C29").Copy()
$ws.Range("C29").PasteSpecial(-4163)
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("D30").Formula = "=""0"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("C30").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").Formula = "=""***.*"""
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("M30").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$ws.Range("G30").Value = 4

$ws.Range("H30").Value = -75
